# template_nonASN.xlsx edit: update I2 date value and move the live selection
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell I2 (date-formatted, numFmt 14) changes from 37174 (2001-10-10) to 36443 (1999-10-10)
$ws.Range("I2").Value = 36443

# Move the saved selection from the old multi-range "whole column" selection
# (I:I, T:T, W:W with active cell W1) to a single active cell L9
$ws.Range("L9").Select()
